# Weekly update: add a new week of "Cilantro" price data at
# Terminal Hortofrutícola Agro Chillán (rows 126:127), pushing the
# existing historical rows down by two rows (old A126:R195 -> A128:R197).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the current row 126; Excel
# shifts rows 126:195 down to 128:197 and extends the used range
# automatically (dimension becomes A1:R197).
$ws.Rows("126:127").Insert()

# --- New row 126 (Primera) ---
$ws.Range("A126").Value = 7
$ws.Range("B126").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C126").Value = 'Ñuble'
$ws.Range("D126").Value = 45001
$ws.Range("E126").Value = 16
$ws.Range("F126").Value = 100112040
$ws.Range("G126").Value = 'Cilantro'
$ws.Range("H126").Value = 'Sin especificar'
$ws.Range("I126").Value = 'Primera'
$ws.Range("J126").Value = 500
$ws.Range("K126").Value = 2000
$ws.Range("L126").Value = 2000
$ws.Range("M126").Value = 2000
$ws.Range("N126").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O126").Value = 'Provincia de Diguillín'
$ws.Range("P126").Value = 2000
$ws.Range("Q126").Value = 1
$ws.Range("R126").Value = 'Hortaliza'

# --- New row 127 (Segunda) ---
$ws.Range("A127").Value = 7
$ws.Range("B127").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C127").Value = 'Ñuble'
$ws.Range("D127").Value = 45001
$ws.Range("E127").Value = 16
$ws.Range("F127").Value = 100112040
$ws.Range("G127").Value = 'Cilantro'
$ws.Range("H127").Value = 'Sin especificar'
$ws.Range("I127").Value = 'Segunda'
$ws.Range("J127").Value = 300
$ws.Range("K127").Value = 1500
$ws.Range("L127").Value = 1500
$ws.Range("M127").Value = 1500
$ws.Range("N127").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O127").Value = 'Provincia de Diguillín'
$ws.Range("P127").Value = 1500
$ws.Range("Q127").Value = 1
$ws.Range("R127").Value = 'Hortaliza'
